$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to literal text so multi-dot values (e.g. '28.968.26')
# and values with significant trailing zeros (e.g. '4.920', '1.000') are stored
# exactly as received from the feed, matching the original inline-string cells.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '28.968.26'
$ws.Range('E2').Value = '  -0.52%  '

$ws.Range('D3').Value = '1.817.30'
$ws.Range('E3').Value = '  -1.05%  '

$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '240.94'
$ws.Range('E5').Value = '  -1.35%  '

$ws.Range('D6').Value = '0.6091'
$ws.Range('E6').Value = '  -3.36%  '

$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.21%  '

$ws.Range('D8').Value = '0.07300'
$ws.Range('E8').Value = '  -2.79%  '

$ws.Range('D9').Value = '0.2873'
$ws.Range('E9').Value = '  -1.95%  '

$ws.Range('D10').Value = '22.81'
$ws.Range('E10').Value = '  -1.93%  '

$ws.Range('D11').Value = '0.07636'
$ws.Range('E11').Value = '  -1.23%  '

$ws.Range('D12').Value = '1.813.59'
$ws.Range('E12').Value = '  -1.18%  '

$ws.Range('D13').Value = '4.920'
$ws.Range('E13').Value = '  -1.59%  '

$ws.Range('D14').Value = '0.6569'
$ws.Range('E14').Value = '  -1.98%  '

$ws.Range('D15').Value = '81.02'
$ws.Range('E15').Value = '  -2.10%  '

$ws.Range('D16').Value = '0.000008871'
$ws.Range('E16').Value = '  -5.01%  '

$ws.Range('D17').Value = '5.830'
$ws.Range('E17').Value = '  -3.01%  '

$ws.Range('D18').Value = '28.938.87'
$ws.Range('E18').Value = '  -0.68%  '

$ws.Range('D19').Value = '2.065.73'
$ws.Range('E19').Value = '  -0.96%  '

$ws.Range('D20').Value = '235.06'
$ws.Range('E20').Value = '  +4.92%  '

$ws.Range('D21').Value = '12.39'
$ws.Range('E21').Value = '  -1.78%  '

$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.24%  '

$ws.Range('D23').Value = '7.110'
$ws.Range('E23').Value = '  -0.53%  '

$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.14%  '

$ws.Range('D25').Value = '158.44'
$ws.Range('E25').Value = '  -0.85%  '

$ws.Range('D26').Value = '0.1394'
$ws.Range('E26').Value = '  -0.76%  '

$ws.Range('D27').Value = '8.387'
$ws.Range('E27').Value = '  -1.53%  '

$ws.Range('D28').Value = '17.54'
$ws.Range('E28').Value = '  -2.49%  '

$ws.Range('D29').Value = '1.476'
$ws.Range('E29').Value = '  -1.66%  '

$ws.Range('D30').Value = '0.05588'
$ws.Range('E30').Value = '  -3.39%  '

$ws.Range('D31').Value = '4.054'
$ws.Range('E31').Value = '  -0.40%  '

$ws.Range('D32').Value = '4.062'
$ws.Range('E32').Value = '  -2.45%  '

$ws.Range('D33').Value = '1.207'
$ws.Range('E33').Value = '  +0.12%  '

$ws.Range('D34').Value = '1.819'
$ws.Range('E34').Value = '  -1.67%  '

$ws.Range('D35').Value = '0.7282'
$ws.Range('E35').Value = '  -3.06%  '

$ws.Range('D36').Value = '1.128'
$ws.Range('E36').Value = '  -1.11%  '

$ws.Range('D37').Value = '2.620'
$ws.Range('E37').Value = '  -2.07%  '

$ws.Range('D38').Value = '2.804'
$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.192.24'
$ws.Range('E39').Value = '  -2.73%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01746'
$ws.Range('E40').Value = '  -2.32%  '

$ws.Range('E41').Value = '  -3.83%  '

$ws.Range('D42').Value = '0.8780'
$ws.Range('E42').Value = '  -1.82%  '

$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  -0.29%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '100.73'
$ws.Range('E44').Value = '  -1.62%  '

$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.962.91'
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('D46').Value = '64.14'
$ws.Range('E46').Value = '  -2.75%  '

$ws.Range('D47').Value = '0.5085'
$ws.Range('E47').Value = '  -0.30%  '

$ws.Range('D48').Value = '0.00000000119'
$ws.Range('E48').Value = '  -3.83%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = '0.3972'
$ws.Range('E49').Value = '  -2.63%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.996'
$ws.Range('E50').Value = '  -0.79%  '

$ws.Range('D51').Value = '0.05780'
$ws.Range('E51').Value = '  -0.93%  '

# Clean up: drop the explicit text-number-format style from the Price cells so
# they keep the workbook's default (un-styled) formatting, same as before the edit.
$ws.Range('D2:D51').Style = 'Normal'